# Add new worksheet "CompleteRFC1" after the last existing sheet,
# mirroring the structure/data of "CompleteRFC" (this is the new
# hyperparameter-results sheet added by the "Hyperparameters" rename commit).

$wb = $excel.ActiveWorkbook

$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "CompleteRFC1"

# --- populate A1:M17 with the RFC hyperparameter results table ---
$data = New-Object 'object[,]' 17,13
$data[0,1] = "Conf_id"
$data[0,2] = "Dataset"
$data[0,3] = "Criterion"
$data[0,4] = "max_depth"
$data[0,5] = "n_estimators"
$data[0,6] = "max_features"
$data[0,7] = "P1"
$data[0,8] = "P2"
$data[0,9] = "P3"
$data[0,10] = "P4"
$data[0,11] = "P5"
$data[0,12] = "Promedio"
$data[1,0] = 0
$data[1,1] = 1
$data[1,2] = "Complete"
$data[1,3] = "gini"
$data[1,4] = 5
$data[1,5] = 10
$data[1,6] = 5
$data[1,7] = 0.8082813012803404
$data[1,8] = 0.8047662376139391
$data[1,9] = 0.7861722212162833
$data[1,10] = 0.7742819099071137
$data[1,11] = 0.8177520061402811
$data[1,12] = 0.7982507352315914
$data[2,0] = 1
$data[2,1] = 2
$data[2,2] = "Complete"
$data[2,3] = "gini"
$data[2,4] = 7
$data[2,5] = 30
$data[2,6] = 7
$data[2,7] = 0.9614932441097477
$data[2,8] = 0.9522535577528214
$data[2,9] = 0.9496338346378987
$data[2,10] = 0.9290591409681976
$data[2,11] = 0.9398414431196858
$data[2,12] = 0.9464562441176703
$data[3,0] = 2
$data[3,1] = 3
$data[3,2] = "Complete"
$data[3,3] = "gini"
$data[3,4] = 10
$data[3,5] = 100
$data[3,6] = 9
$data[3,7] = 0.985740541262185
$data[3,8] = 0.9827410729132375
$data[3,9] = 0.985667235573898
$data[3,10] = 0.9855146623457041
$data[3,11] = 0.9855509657450032
$data[3,12] = 0.9850428955680055
$data[4,0] = 3
$data[4,1] = 4
$data[4,2] = "Complete"
$data[4,3] = "gini"
$data[4,4] = 13
$data[4,5] = 150
$data[4,6] = 15
$data[4,7] = 0.9852105060606411
$data[4,8] = 0.9841173026983225
$data[4,9] = 0.989763917212975
$data[4,10] = 0.9890683669299062
$data[4,11] = 0.9906418203275371
$data[4,12] = 0.9877603826458763
$data[5,0] = 4
$data[5,1] = 5
$data[5,2] = "Complete"
$data[5,3] = "gini"
$data[5,4] = 15
$data[5,5] = 200
$data[5,6] = 21
$data[5,7] = 0.9871561711599128
$data[5,8] = 0.9848347327925404
$data[5,9] = 0.9905007847468091
$data[5,10] = 0.9863879619576181
$data[5,11] = 0.9917160574374575
$data[5,12] = 0.9881191416188676
$data[6,0] = 5
$data[6,1] = 6
$data[6,2] = "Complete"
$data[6,3] = "gini"
$data[6,4] = 17
$data[6,5] = 400
$data[6,6] = 19
$data[6,7] = 0.9854234449234226
$data[6,8] = 0.9861011768307386
$data[6,9] = 0.9889985263520042
$data[6,10] = 0.9864200181289061
$data[6,11] = 0.9907442279034873
$data[6,12] = 0.9875374788277117
$data[7,0] = 6
$data[7,1] = 7
$data[7,2] = "Complete"
$data[7,3] = "gini"
$data[7,4] = 21
$data[7,5] = 500
$data[7,6] = 18
$data[7,7] = 0.9846924928332484
$data[7,8] = 0.985146133179174
$data[7,9] = 0.9889985263520042
$data[7,10] = 0.9878342470258968
$data[7,11] = 0.9917160574374575
$data[7,12] = 0.9876774913655563
$data[8,0] = 7
$data[8,1] = 8
$data[8,2] = "Complete"
$data[8,3] = "gini"
$data[8,4] = 10
$data[8,5] = 150
$data[8,6] = 10
$data[8,7] = 0.9847755178524658
$data[8,8] = 0.9814632114151935
$data[8,9] = 0.9827844301601295
$data[8,10] = 0.9867877909273719
$data[8,11] = 0.9862817597689497
$data[8,12] = 0.9844185420248222
$data[9,0] = 8
$data[9,1] = 9
$data[9,2] = "Complete"
$data[9,3] = "entropy"
$data[9,4] = 5
$data[9,5] = 10
$data[9,6] = 13
$data[9,7] = 0.8268748805871683
$data[9,8] = 0.8086577251796926
$data[9,9] = 0.842807864571653
$data[9,10] = 0.8154401136565022
$data[9,11] = 0.8428246127280945
$data[9,12] = 0.8273210393446222
$data[10,0] = 9
$data[10,1] = 10
$data[10,2] = "Complete"
$data[10,3] = "entropy"
$data[10,4] = 7
$data[10,5] = 30
$data[10,6] = 15
$data[10,7] = 0.9624350795053922
$data[10,8] = 0.9455367976390423
$data[10,9] = 0.9682938902993993
$data[10,10] = 0.9624472419109688
$data[10,11] = 0.9650367735341482
$data[10,12] = 0.9607499565777902
$data[11,0] = 10
$data[11,1] = 11
$data[11,2] = "Complete"
$data[11,3] = "entropy"
$data[11,4] = 10
$data[11,5] = 100
$data[11,6] = 7
$data[11,7] = 0.9862806813923686
$data[11,8] = 0.9820324469646149
$data[11,9] = 0.9802179160167925
$data[11,10] = 0.9849648832414222
$data[11,11] = 0.9905162072216733
$data[11,12] = 0.9848024269673742
$data[12,0] = 11
$data[12,1] = 12
$data[12,2] = "Complete"
$data[12,3] = "entropy"
$data[12,4] = 13
$data[12,5] = 150
$data[12,6] = 3
$data[12,7] = 0.9777663503945369
$data[12,8] = 0.9795182192488896
$data[12,9] = 0.9723756786386794
$data[12,10] = 0.9843260474196867
$data[12,11] = 0.9826680716890173
$data[12,12] = 0.9793308734781618
$data[13,0] = 12
$data[13,1] = 13
$data[13,2] = "Complete"
$data[13,3] = "entropy"
$data[13,4] = 15
$data[13,5] = 200
$data[13,6] = 21
$data[13,7] = 0.9850344774743822
$data[13,8] = 0.9793058050422392
$data[13,9] = 0.9891717250110007
$data[13,10] = 0.9849878298975374
$data[13,11] = 0.9913037542948692
$data[13,12] = 0.9859607183440058
$data[14,0] = 13
$data[14,1] = 14
$data[14,2] = "Complete"
$data[14,3] = "entropy"
$data[14,4] = 17
$data[14,5] = 400
$data[14,6] = 11
$data[14,7] = 0.9858502955861703
$data[14,8] = 0.9842008947337588
$data[14,9] = 0.9891511721149359
$data[14,10] = 0.9882685347127076
$data[14,11] = 0.9918370519325963
$data[14,12] = 0.9878615898160337
$data[15,0] = 14
$data[15,1] = 15
$data[15,2] = "Complete"
$data[15,3] = "entropy"
$data[15,4] = 21
$data[15,5] = 500
$data[15,6] = 15
$data[15,7] = 0.9859419079025691
$data[15,8] = 0.9836645759570886
$data[15,9] = 0.9889375881764481
$data[15,10] = 0.9868532508086646
$data[15,11] = 0.993247526687467
$data[15,12] = 0.9877289699064475
$data[16,0] = 15
$data[16,1] = 16
$data[16,2] = "Complete"
$data[16,3] = "entropy"
$data[16,4] = 10
$data[16,5] = 150
$data[16,6] = 12
$data[16,7] = 0.9852143264489701
$data[16,8] = 0.9819171589060189
$data[16,9] = 0.9833446074080839
$data[16,10] = 0.9871951616772192
$data[16,11] = 0.9870642026677869
$data[16,12] = 0.9849470914216157

$ws.Range("A1:M17").Value = $data

# --- formatting: bold + centered/top aligned + thin-bordered header row
#     (B1:M1) and "index" column (A2:A17), matching the look of the other
#     sheets in this workbook (CompleteRFC, CompleteNB1, ...) ---
$headerRange = $ws.Range("B1:M1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

$indexRange = $ws.Range("A2:A17")
$indexRange.Font.Bold = $true
$indexRange.HorizontalAlignment = -4108
$indexRange.VerticalAlignment = -4160
$indexRange.Borders.LineStyle = 1

$ws.Range("A1").Select()

Write-Host "Added sheet 'CompleteRFC1' with $($ws.UsedRange.Rows.Count) rows and $($ws.UsedRange.Columns.Count) columns"
